$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
#    (The paragraph containing the bold "Meta description" run plus the
#    description text run.)
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Buffalo Power Hold and Win for free"
#    right before the very last paragraph (the one that currently holds the
#    italic "Prompt: ..." text).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Buffalo Power Hold and Win for free</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlFrag)

# InsertXML leaves behind an extra empty paragraph used purely as a
# separator; remove it so the "Prompt"/new-text paragraph immediately
# follows the newly inserted bold paragraph.
$extraEmptyIndex = $d.Paragraphs.Count - 1
$extraEmptyPara = $d.Paragraphs.Item($extraEmptyIndex)
$extraEmptyPara.Range.Delete()

# 3. Replace the old "Prompt: ..." italic text with the new review summary
#    text, keeping the existing italic run formatting untouched.
$d.Content.Find.Execute(
    "Prompt: Create a cartoon style image featuring a happy Maya warrior with glasses. Instructions: DALLE, please use cartoon-style graphics and create an image of a happy Maya warrior with glasses. The warrior should be wearing a traditional Mayan headdress and outfit, and the glasses should look modern. The image should have bright colors and should be eye-catching to attract players to the game " + [char]34 + "Buffalo Power Hold and Win" + [char]34 + ". Please make sure the image complements the game's theme of wild North America and the buffalo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Buffalo Power Hold and Win and play for free. Features, Jackpots, Graphics, Betting options and RTP.",
    2)
